$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of FX data for 2021-03-31 (USD/HKD), matching the format
# of the existing rows above it (row 8).
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A9").Value = (Get-Date -Year 2021 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B9").Value = "USD"
$ws.Range("C9").Value = "HKD"
$ws.Range("D9").Value = 7.7740999999999998

$ws.Range("D10").Select()
